$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.958.09'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '2.548.60'
$ws.Range('E3').Value = '  +2.98%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'568.99"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').Value = "'146.31"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.99%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = "'0.582"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = '2.547.14'
$ws.Range('E9').Value = '  +2.94%  '
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('E11').Value = '  -4.18%  '
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').Value = '3.003.08'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').Value = '62.894.19'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '2.543.35'
$ws.Range('E18').Value = '  +3.07%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('D21').Value = "'334.61"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('D22').Value = "'6.76"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = "'65.38"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('E25').Value = '  -1.30%  '
$ws.Range('E26').Value = '  +4.98%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('E28').Value = '  +2.64%  '
$ws.Range('D29').Value = "'8.33"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('D30').Value = "'7.34"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.46%  '
$ws.Range('D31').Value = '0.0₃0812'
$ws.Range('E31').Value = '  +1.37%  '
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').Value = "'175.42"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('D35').Value = "'404.78"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').Value = "'0.400"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = "'19.09"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = "'39.43"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('D43').Value = "'151.73"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.04%  '
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').Value = "'20.71"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = "'0.0530"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.24%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = "'0.0965"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('E49').Value = '  +4.60%  '
$ws.Range('E50').Value = '  +0.96%  '
$ws.Range('E51').Value = '  -3.04%  '
